# "further cleaning to metadata"
#  1. Shared string "E7760" -> "E7420" (the sample/run number shown in column G)
#  2. Sheet view: scroll position back to the top (A1) and move the active
#     selection from column H to column G (H2:H41 -> G2:G41)
#  3. H2:H41: turn the literal boolean FALSE into a live =FALSE() formula

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the shared "E7760" sample id text to "E7420" -----------------
# All of G2:G41 point at the same shared-string entry, so rewrite them all in
# one batch so the old text is fully replaced (no leftover/duplicate string).
$ws.Range("G2:G41").Value = "E7420"

# --- 3. H2:H41 -> live formulas instead of cached boolean literals ----------
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=FALSE()"
}

# --- 2. Restore the gridlines flag, scroll back to the top, move selection --
$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$win.ScrollRow = 1
$win.ScrollColumn = 1

$ws.Range("G2:G41").Select() | Out-Null
